# SIRA.xlsx edit ("Add files via upload"):
#   - the "Gráfico1" chart sheet is dropped, leaving the data worksheet as
#     the only (and therefore active) sheet
#   - the surviving worksheet "Respuestas de formulario 1" is renamed to
#     "DATOS"
#   - the live selection on that sheet moves to D14

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Rename the surviving data worksheet and make it the active one, since the
# chart sheet that used to hold the active tab is going away.
$ws = $wb.Worksheets.Item("Respuestas de formulario 1")
$ws.Name = "DATOS"
$ws.Activate()

# Leave the selection on D14 (still inside the frozen "bottomLeft" pane).
$ws.Range("D14").Select()

# Best-effort: remove the now-unused "Gráfico1" chart sheet if the host
# exposes it through the Sheets collection. Some hosts only surface regular
# worksheets there, in which case this is a harmless no-op.
try {
    $wb.Sheets.Item("Gráfico1").Delete()
} catch {
}
